# "Models clean up completed" - rebuild the sales report rows.
# Row 2 and 3 get corrected (date/qty/product/price/total), and the
# sheet grows from 3 data rows to 15 data rows (rows 2-16), most of
# which no longer have a product name (cleared out).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order Date, Order ID, Product name, Price, Total Amount, Payment method, Payment Status
$rows = @(
    @("02/24/2023", 1,  "",        0,     180000, "COD",       "pending"),
    @("02/24/2023", 2,  "",        0,     120000, "Razor Pay", "pending"),
    @("02/24/2023", 3,  "",        0,     120000, "Wallet",    "pending"),
    @("02/24/2023", 4,  "",        0,     300000, "COD",       "pending"),
    @("02/24/2023", 5,  "Oppo x3", 60000, 60000,  "Wallet",    "pending"),
    @("02/24/2023", 6,  "",        0,     60000,  "COD",       "pending"),
    @("02/24/2023", 7,  "",        0,     60000,  "COD",       "pending"),
    @("02/24/2023", 8,  "",        0,     120000, "COD",       "pending"),
    @("02/24/2023", 9,  "",        0,     60000,  "COD",       "pending"),
    @("02/24/2023", 10, "",        0,     60000,  "COD",       "pending"),
    @("02/24/2023", 11, "",        0,     60000,  "COD",       "pending"),
    @("02/24/2023", 13, "",        0,     60000,  "Wallet",    "pending"),
    @("02/24/2023", 12, "",        0,     60000,  "Wallet",    "pending"),
    @("02/24/2023", 14, "",        0,     0,      "",          ""),
    @("02/24/2023", 15, "",        0,     0,      "",          "")
)

$r = 2
foreach ($row in $rows) {
    # Force the date column to stay plain text (matches source file,
    # where dates are literal strings, not date serials).
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $r = $r + 1
}
